$d = $word.ActiveDocument
$d.Content.Find.Execute("Le Thanh", $true, $false, $false, $false, $false, $true, 1, $false, "Thanh Le", 2)
